# Weekly update: a new price record is inserted at the top of the data
# (row 2). All existing data rows 2-29 shift down by one row (to rows
# 3-30), so the worksheet grows from 29 data rows to 30.
#
# Columns A, B, C, E, F, G, H, I, N, Q, R are identical for every data
# row in this sheet, so when a row shifts down we just need to make
# sure those are present too; easiest is to copy the whole row content.
# Only D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) differ
# between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 29
$newLastDataRow = 30

# 1. Read every existing data row (2..29) into memory before we start
#    overwriting anything.
$savedRows = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowValues = @{}
    for ($col = 1; $col -le 18; $col++) {
        $rowValues[$col] = $ws.Cells.Item($r, $col).Value2
    }
    $savedRows[$r] = $rowValues
}

# 2. Write the saved rows back out, shifted down by one (old row r ->
#    new row r+1), working from the bottom up so we never clobber data
#    we still need to copy.
$dateFormat = $ws.Cells.Item($firstDataRow, 4).NumberFormat
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + 1
    $rowValues = $savedRows[$r]
    for ($col = 1; $col -le 18; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $rowValues[$col]
    }
    # Column D (Fecha) keeps its date number format on every data row,
    # including the brand new last row (30) which has no pre-existing
    # style to inherit from.
    $ws.Cells.Item($destRow, 4).NumberFormat = $dateFormat
}

# 3. Populate the new first data row (row 2) with this week's record.
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(2, 3).Value = "Maule"
$ws.Cells.Item(2, 4).Value = 44750
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 100112040
$ws.Cells.Item(2, 7).Value = "Cilantro"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 150
$ws.Cells.Item(2, 11).Value = 9000
$ws.Cells.Item(2, 12).Value = 9000
$ws.Cells.Item(2, 13).Value = 9000
$ws.Cells.Item(2, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(2, 15).Value = "Región Metropolitana"
$ws.Cells.Item(2, 16).Value = 250
$ws.Cells.Item(2, 17).Value = 36
$ws.Cells.Item(2, 18).Value = "Hortaliza"
$ws.Cells.Item(2, 4).NumberFormat = $dateFormat

